# Populate the "empty" 3-sheet workbook with sample data, give the sheets
# descriptive names, and leave the selection/active-tab state the way the
# authors left it when they saved the file.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Rename the sheets -----------------------------------------------------
$ws1.Name = "Sheet1 - Text"
$ws2.Name = "Sheet2 - Numbers"
$ws3.Name = "Sheet3 - Formulas"

# --- Sheet1 - Text: a couple of text cells ---------------------------------
$ws1.Range("A1").Value = "This is cell A1 in Sheet 1"
$ws1.Range("G5").Value = "This is cell G5"

# --- Sheet2 - Numbers: a column of numbers + a column of percentages ------
for ($row = 1; $row -le 30; $row++) {
    $ws2.Cells.Item($row, 4).Value = $row           # D1:D30 -> 1..30
    $ws2.Cells.Item($row, 11).Value = $row / 100     # K1:K30 -> 1%..30%
    $ws2.Cells.Item($row, 11).Style = "Percent"
}
$ws2.Range("G5").Value = "This is cell G5"

# --- Sheet3 - Formulas: a formula referencing Sheet2 -----------------------
$ws3.Range("D2").Formula = "='Sheet2 - Numbers'!D5"

# --- Selections / active sheet, in the order they were left in -------------
$ws1.Range("G6").Select()
$ws2.Range("L2").Select()
$ws3.Range("D3").Select()
